$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 71; this shifts the existing rows 71..148
# down to 72..149 and extends the sheet's used range accordingly.
$ws.Rows(71).Insert()

# Populate the newly inserted row 71 with the new record.
$ws.Cells.Item(71, 1).Value = 10
$ws.Cells.Item(71, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(71, 3).Value = "La Araucanía"
$ws.Cells.Item(71, 4).Value = 44778
$ws.Cells.Item(71, 5).Value = 9
$ws.Cells.Item(71, 6).Value = 100112031
$ws.Cells.Item(71, 7).Value = "Poroto verde"
$ws.Cells.Item(71, 8).Value = "Sin especificar"
$ws.Cells.Item(71, 9).Value = "Primera"
$ws.Cells.Item(71, 10).Value = 65
$ws.Cells.Item(71, 11).Value = 28000
$ws.Cells.Item(71, 12).Value = 30000
$ws.Cells.Item(71, 13).Value = 28923
$ws.Cells.Item(71, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(71, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(71, 16).Value = 1157
$ws.Cells.Item(71, 17).Value = 25
$ws.Cells.Item(71, 18).Value = "Hortaliza"
